# Rapport Unity Eksamen.docx — apply the commit's content changes.
#
# The vast majority of the unified diff is Word's automatic proofing-tool
# markup (<w:proofErr w:type="spellStart"/>.../<w:proofErr .../> pairs plus
# the accompanying run splits) around words Word's Norwegian dictionary
# doesn't recognise (Unity, Github, Git, merge, backup, features, bugs,
# Kitti's, etc.). That markup is produced internally by Word's spell
# checker as a side effect of proofing the document and is not something
# the Word object model exposes a way to author directly - it carries no
# visible/textual content of its own. Every run that diff splits joins
# back into exactly the same paragraph text, so no actual wording changes
# there.
#
# The two real, content-visible edits in the commit are:
#   1. A new paragraph from Kaja ("Kaja- Jeg slet med ...") is inserted
#      right after Murvet's paragraph in section 3.
#   2. One of the two blank paragraphs that used to sit between Murvet's
#      paragraph and the "4. Kjente bugs" heading is removed (two blanks
#      become one).
#
# Both are reproduced below using the Word object model.

$d = $word.ActiveDocument

# Anchor on the end of Murvet's paragraph text (unique in the document).
$anchorRange = $d.Content
$anchorRange.Find.Execute("jobba hardt for.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$murvetPara = $anchorRange.Paragraphs(1)
$firstBlank = $murvetPara.Next()

# Insert a brand-new, unformatted paragraph before the first blank one,
# then fill it with Kaja's text (this keeps the new run free of any
# inherited character formatting, matching the target markup).
$firstBlank.Range.InsertParagraphBefore()
$kajaPara = $murvetPara.Next()
$kajaPara.Range.Text = "Kaja- Jeg slet med å få triggerEnter koden til å fungere, så på utallige youtube tutorials og til slutt hjalp Kitty meg med en kode han hadde, som jeg bare endret til min egen. Ellers var det en del problemer med selve prosjektet som skrevet ovenfor. Hadde litt problemer med menyene som ikke alltid ville funke på alle scenene. "

# Collapse the two blank paragraphs that followed Murvet's paragraph down
# to one by deleting the first of them (now immediately after Kaja's
# paragraph).
$blankAfterKaja = $kajaPara.Next()
$blankAfterKaja.Range.Delete() | Out-Null
